# Auto-generated script applying scheduled-runner market-price refresh to Zodiark_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 300
$ws.Range("K5").Value = 300
$ws.Range("M5").Value = -185
$ws.Range("H17").Value = 564.0323
$ws.Range("J17").Value = 564.0323
$ws.Range("L17").Value = 1692.0969
$ws.Range("N17").Value = -2028.0969
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H52").Value = 1546.909
$ws.Range("I52").Value = 501.8889
$ws.Range("J52").Value = 6249.5
$ws.Range("K52").Value = 1505.6667
$ws.Range("L52").Value = 18748.5
$ws.Range("M52").Value = -1345.6667
$ws.Range("N52").Value = -19068.5
$ws.Range("H86").Value = 1530.8
$ws.Range("I86").Value = 1600
$ws.Range("K86").Value = 1600
$ws.Range("M86").Value = -477
$ws.Range("H89").Value = 1530.8
$ws.Range("I89").Value = 1600
$ws.Range("K89").Value = 8000
$ws.Range("M89").Value = -2384
$ws.Range("H125").Value = 1499.6666
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1499.6666
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 13496.9994
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -18416.9994
$ws.Range("H132").Value = 2187.6086
$ws.Range("I132").Value = 2070.75
$ws.Range("K132").Value = 6212.25
$ws.Range("M132").Value = -3682.25
$ws.Range("H137").Value = 2377.5715
$ws.Range("I137").Value = 2630.7334
$ws.Range("J137").Value = 1744.6666
$ws.Range("K137").Value = 7892.2002
$ws.Range("L137").Value = 5233.9998
$ws.Range("M137").Value = -5342.2002
$ws.Range("N137").Value = -10333.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5146.7334
$ws.Range("I61").Value = 6574.75
$ws.Range("J61").Value = 4627.4546
$ws.Range("K61").Value = 6574.75
$ws.Range("L61").Value = 4627.4546
$ws.Range("M61").Value = -6362.75
$ws.Range("N61").Value = -5051.4546
$ws.Range("H74").Value = 4291.8945
$ws.Range("I74").Value = 4007.625
$ws.Range("J74").Value = 4498.636
$ws.Range("K74").Value = 4007.625
$ws.Range("L74").Value = 4498.636
$ws.Range("M74").Value = -3133.625
$ws.Range("N74").Value = -6246.636
$ws.Range("H77").Value = 4291.8945
$ws.Range("I77").Value = 4007.625
$ws.Range("J77").Value = 4498.636
$ws.Range("K77").Value = 20038.125
$ws.Range("L77").Value = 22493.18
$ws.Range("M77").Value = -15670.125
$ws.Range("N77").Value = -31229.18
$ws.Range("H80").Value = 50110
$ws.Range("J80").Value = 50110
$ws.Range("L80").Value = 50110
$ws.Range("N80").Value = -52106
$ws.Range("H83").Value = 50110
$ws.Range("J83").Value = 50110
$ws.Range("L83").Value = 150330
$ws.Range("N83").Value = -160314
$ws.Range("H88").Value = 2079.3845
$ws.Range("J88").Value = 2167.5293
$ws.Range("L88").Value = 2167.5293
$ws.Range("N88").Value = -2979.5293
$ws.Range("H91").Value = 2079.3845
$ws.Range("J91").Value = 2167.5293
$ws.Range("L91").Value = 2167.5293
$ws.Range("N91").Value = -4975.5293
$ws.Range("H110").Value = 1827.6666
$ws.Range("I110").Value = 921.2857
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 921.2857
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = 1123.7143
$ws.Range("N110").Value = -9090
$ws.Range("H122").Value = 4038.8647
$ws.Range("I122").Value = 4187.885
$ws.Range("J122").Value = 3686.6365
$ws.Range("K122").Value = 12563.655
$ws.Range("L122").Value = 11059.9095
$ws.Range("M122").Value = -10113.655
$ws.Range("N122").Value = -15959.9095
$ws.Range("H136").Value = 5146.7334
$ws.Range("I136").Value = 6574.75
$ws.Range("J136").Value = 4627.4546
$ws.Range("K136").Value = 19724.25
$ws.Range("L136").Value = 13882.3638
$ws.Range("M136").Value = -17174.25
$ws.Range("N136").Value = -18982.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1583.75
$ws.Range("I20").Value = 1534.6666
$ws.Range("J20").Value = 1731
$ws.Range("K20").Value = 1534.6666
$ws.Range("L20").Value = 1731
$ws.Range("M20").Value = -1287.6666
$ws.Range("N20").Value = -2225
$ws.Range("H107").Value = 2307.0667
$ws.Range("I107").Value = 1284.8636
$ws.Range("J107").Value = 5118.125
$ws.Range("K107").Value = 1284.8636
$ws.Range("L107").Value = 5118.125
$ws.Range("M107").Value = 635.1364000000001
$ws.Range("N107").Value = -8958.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 11820
$ws.Range("I105").Value = 13149.4
$ws.Range("K105").Value = 13149.4
$ws.Range("M105").Value = -11402.4
$ws.Range("H107").Value = 1204.591
$ws.Range("J107").Value = 1213.8096
$ws.Range("L107").Value = 1213.8096
$ws.Range("N107").Value = -5053.809600000001
$ws.Range("H134").Value = 2595.4614
$ws.Range("I134").Value = 2491.25
$ws.Range("K134").Value = 7473.75
$ws.Range("M134").Value = -4938.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1626.3684
$ws.Range("I113").Value = 1493.4445
$ws.Range("J113").Value = 1746
$ws.Range("K113").Value = 4480.333500000001
$ws.Range("L113").Value = 5238
$ws.Range("M113").Value = -2310.333500000001
$ws.Range("N113").Value = -9578
$ws.Range("H114").Value = 11960.632
$ws.Range("I114").Value = 468.42856
$ws.Range("J114").Value = 18664.416
$ws.Range("K114").Value = 1405.28568
$ws.Range("L114").Value = 55993.24800000001
$ws.Range("M114").Value = 1848.71432
$ws.Range("N114").Value = -62501.24800000001
$ws.Range("H120").Value = 10711.8
$ws.Range("I120").Value = 5423.6
$ws.Range("K120").Value = 16270.8
$ws.Range("M120").Value = -11432.8
$ws.Range("H128").Value = 680996.25
$ws.Range("I128").Value = 680996.25
$ws.Range("K128").Value = 2042988.75
$ws.Range("M128").Value = -2038008.75
$ws.Range("H129").Value = 3826.842
$ws.Range("I129").Value = 4899.125
$ws.Range("K129").Value = 14697.375
$ws.Range("M129").Value = -9697.375
$ws.Range("H131").Value = 703.5
$ws.Range("I131").Value = 561.1429000000001
$ws.Range("J131").Value = 1700
$ws.Range("K131").Value = 1683.4287
$ws.Range("L131").Value = 5100
$ws.Range("M131").Value = 3356.5713
$ws.Range("N131").Value = -15180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 161.53847
$ws.Range("I2").Value = 181.875
$ws.Range("K2").Value = 181.875
$ws.Range("M2").Value = -68.875
$ws.Range("H7").Value = 6666
$ws.Range("I7").Value = 6666
$ws.Range("K7").Value = 6666
$ws.Range("M7").Value = -6554
$ws.Range("H8").Value = 6666
$ws.Range("I8").Value = 6666
$ws.Range("K8").Value = 6666
$ws.Range("M8").Value = -6527
$ws.Range("H15").Value = 24999.5
$ws.Range("J15").Value = 24999.5
$ws.Range("L15").Value = 24999.5
$ws.Range("N15").Value = -25575.5
$ws.Range("H70").Value = 24221.645
$ws.Range("I70").Value = 36112.35
$ws.Range("J70").Value = 9782.929
$ws.Range("K70").Value = 36112.35
$ws.Range("L70").Value = 9782.929
$ws.Range("M70").Value = -35842.35
$ws.Range("N70").Value = -10322.929
$ws.Range("H73").Value = 24221.645
$ws.Range("I73").Value = 36112.35
$ws.Range("J73").Value = 9782.929
$ws.Range("K73").Value = 36112.35
$ws.Range("L73").Value = 9782.929
$ws.Range("M73").Value = -35176.35
$ws.Range("N73").Value = -11654.929
$ws.Range("H81").Value = 24999.5
$ws.Range("J81").Value = 24999.5
$ws.Range("L81").Value = 24999.5
$ws.Range("N81").Value = -26995.5
$ws.Range("H84").Value = 24999.5
$ws.Range("J84").Value = 24999.5
$ws.Range("L84").Value = 74998.5
$ws.Range("N84").Value = -84982.5
$ws.Range("H98").Value = 37400
$ws.Range("J98").Value = 37400
$ws.Range("L98").Value = 37400
$ws.Range("N98").Value = -43390
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H132").Value = 2001.5454
$ws.Range("I132").Value = 1735.2222
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 5205.6666
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -2675.6666
$ws.Range("N132").Value = -14660

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 3999.5
$ws.Range("I107").Value = 3999.5
$ws.Range("K107").Value = 3999.5
$ws.Range("M107").Value = -2079.5
$ws.Range("H132").Value = 3902.3333
$ws.Range("I132").Value = 3421.0833
$ws.Range("J132").Value = 5827.3335
$ws.Range("K132").Value = 10263.2499
$ws.Range("L132").Value = 17482.0005
$ws.Range("M132").Value = -7733.249899999999
$ws.Range("N132").Value = -22542.0005
$ws.Range("H136").Value = 3662.4119
$ws.Range("I136").Value = 3272.25
$ws.Range("K136").Value = 9816.75
$ws.Range("M136").Value = -7266.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 19425
$ws.Range("I3").Value = 19850
$ws.Range("J3").Value = 19000
$ws.Range("K3").Value = 19850
$ws.Range("L3").Value = 19000
$ws.Range("M3").Value = -19736
$ws.Range("N3").Value = -19228
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("H48").Value = 7500
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H132").Value = 998.4286
$ws.Range("I132").Value = 928.8
$ws.Range("J132").Value = 1172.5
$ws.Range("K132").Value = 2786.4
$ws.Range("L132").Value = 3517.5
$ws.Range("M132").Value = -256.3999999999996
$ws.Range("N132").Value = -8577.5
$ws.Range("H136").Value = 4865.4326
$ws.Range("I136").Value = 4458.129
$ws.Range("J136").Value = 6969.8335
$ws.Range("K136").Value = 13374.387
$ws.Range("L136").Value = 20909.5005
$ws.Range("M136").Value = -10824.387
$ws.Range("N136").Value = -26009.5005
